$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.369.55"
$ws.Cells.Item(2, 5).Value = "  +5.68%  "
$ws.Cells.Item(3, 4).Value = "3.465.07"
$ws.Cells.Item(3, 5).Value = "  +6.79%  "
$ws.Cells.Item(4, 5).Value = "  -0.05%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "583.10"
$ws.Cells.Item(5, 5).Value = "  +7.18%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "158.58"
$ws.Cells.Item(6, 5).Value = "  +8.09%  "
$ws.Cells.Item(7, 5).Value = "  -0.14%  "
$ws.Cells.Item(8, 4).Value = "3.469.35"
$ws.Cells.Item(8, 5).Value = "  +6.53%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.535"
$ws.Cells.Item(9, 5).Value = "  +1.41%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "7.61"
$ws.Cells.Item(10, 5).Value = "  +3.15%  "
$ws.Cells.Item(11, 5).Value = "  +7.49%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.441"
$ws.Cells.Item(12, 5).Value = "  +2.04%  "
$ws.Cells.Item(13, 4).Value = "4.058.63"
$ws.Cells.Item(13, 5).Value = "  +6.53%  "
$ws.Cells.Item(14, 5).Value = "  -0.30%  "
$ws.Cells.Item(15, 5).Value = "  +7.81%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "27.68"
$ws.Cells.Item(16, 5).Value = "  +5.25%  "
$ws.Cells.Item(17, 4).Value = "64.376.69"
$ws.Cells.Item(17, 5).Value = "  +5.68%  "
$ws.Cells.Item(18, 4).Value = "3.458.82"
$ws.Cells.Item(18, 5).Value = "  +6.22%  "
$ws.Cells.Item(19, 5).Value = "  +2.57%  "
$ws.Cells.Item(20, 5).Value = "  +7.78%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "398.02"
$ws.Cells.Item(21, 5).Value = "  +5.61%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "8.55"
$ws.Cells.Item(22, 5).Value = "  +1.85%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.547"
$ws.Cells.Item(23, 5).Value = "  +3.23%  "
$ws.Cells.Item(24, 5).Value = "  +0.07%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "72.28"
$ws.Cells.Item(25, 5).Value = "  +3.29%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.0000111"
$ws.Cells.Item(26, 5).Value = "  +21.18%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.51"
$ws.Cells.Item(27, 5).Value = "  +10.32%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.182"
$ws.Cells.Item(28, 5).Value = "  +6.39%  "
$ws.Cells.Item(29, 5).Value = "  -0.30%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.41"
$ws.Cells.Item(30, 5).Value = "  +14.78%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "5.95"
$ws.Cells.Item(31, 5).Value = "  +9.93%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.73"
$ws.Cells.Item(32, 5).Value = "  +9.03%  "
$ws.Cells.Item(33, 5).Value = "  +6.21%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "23.85"
$ws.Cells.Item(34, 5).Value = "  +5.71%  "
$ws.Cells.Item(35, 5).Value = "  -0.01%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.89"
$ws.Cells.Item(36, 5).Value = "  +3.81%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.51"
$ws.Cells.Item(37, 5).Value = "  +5.75%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "158.83"
$ws.Cells.Item(38, 5).Value = "  -0.10%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "28.52"
$ws.Cells.Item(39, 5).Value = "  +7.92%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0787"
$ws.Cells.Item(40, 5).Value = "  +9.32%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.89"
$ws.Cells.Item(41, 5).Value = "  +9.83%  "
$ws.Cells.Item(42, 4).Value = "2.872.57"
$ws.Cells.Item(42, 5).Value = "  +2.35%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0324"
$ws.Cells.Item(43, 5).Value = "  +3.12%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.784"
$ws.Cells.Item(44, 5).Value = "  +7.09%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "42.22"
$ws.Cells.Item(45, 5).Value = "  +5.43%  "
$ws.Cells.Item(46, 5).Value = "  +3.19%  "
$ws.Cells.Item(47, 5).Value = "  +10.23%  "
$ws.Cells.Item(48, 4).Value = "3.507.01"
$ws.Cells.Item(48, 5).Value = "  +6.53%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "22.84"
$ws.Cells.Item(49, 5).Value = "  +6.84%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "300.36"
$ws.Cells.Item(50, 5).Value = "  +8.95%  "
$ws.Cells.Item(51, 5).Value = "  +23.22%  "
